# Refresh the cryptos price/volume(1h) table (GitHub Actions data-refresh commit).
#
# Column D ("Price") values are stored as literal text (e.g. "30.637.07",
# "0.9990") rather than numbers, so several of them would otherwise get
# mangled by Excel's automatic "looks like a number" coercion on plain
# `.Value =` assignment (trailing zeros dropped, multi-dot thousand
# groupings rejected, etc). Writing through `.Formula` with a leading
# apostrophe is the standard COM/Excel "force text" quote-prefix, after
# which `.ClearFormats()` drops the transient quote-prefix cell style so
# the cell ends up as a plain text value with no extra formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'30.637.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Formula = "'1.901.55"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("D4").Formula = "'0.9990"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").Formula = "'238.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Formula = "'0.9988"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Formula = "'0.4783"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Formula = "'0.2836"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Formula = "'0.06539"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Formula = "'2.002.62"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +7.39%  "
$ws.Range("D11").Formula = "'0.07483"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Formula = "'16.67"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Formula = "'5.102"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Formula = "'88.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Formula = "'0.6666"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Formula = "'30.624.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Formula = "'13.32"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.58%  "
# Row 18: coin re-ranked in
$ws.Range("B18").Formula = "'Dai"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").Formula = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").Formula = "'0.9996"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.62%  "

# Row 19: coin re-ranked in
$ws.Range("B19").Formula = "'ShibaInu"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Formula = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").Formula = "'0.000007594"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.18%  "

# Row 20: coin re-ranked in
$ws.Range("B20").Formula = "'WrappedliquidstakedEther2.0"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Formula = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C20").ClearFormats()
$ws.Range("D20").Formula = "'2.191.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.93%  "

$ws.Range("D21").Formula = "'228.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.44%  "
$ws.Range("D22").Formula = "'5.311"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Formula = "'0.9987"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Formula = "'6.220"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Formula = "'168.77"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("D26").Formula = "'9.301"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Formula = "'18.54"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Formula = "'1.957"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("D29").Formula = "'1.402"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").Formula = "'0.09745"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.43%  "
$ws.Range("D31").Formula = "'4.355"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("D32").Formula = "'4.022"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Formula = "'0.05065"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").Formula = "'1.234"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.86%  "
$ws.Range("D35").Formula = "'0.7557"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Formula = "'0.01879"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("D38").Formula = "'2.666"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.70%  "
# Row 39: coin re-ranked in
$ws.Range("B39").Formula = "'RenderToken"
$ws.Range("B39").ClearFormats()
$ws.Range("C39").Formula = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C39").ClearFormats()
$ws.Range("D39").Formula = "'2.087"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.18%  "

# Row 40: coin re-ranked in
$ws.Range("B40").Formula = "'TrustWalletToken"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Formula = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C40").ClearFormats()
$ws.Range("D40").Formula = "'0.9156"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("D41").Formula = "'106.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Formula = "'0.4292"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Formula = "'5.800"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Formula = "'1.005"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Formula = "'7.423"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Formula = "'64.62"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Formula = "'0.1272"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.76%  "
# Row 48: coin re-ranked in
$ws.Range("B48").Formula = "'NEARProtocol"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Formula = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Formula = "'1.479"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.65%  "

# Row 49: coin re-ranked in
$ws.Range("B49").Formula = "'EnergySwap"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Formula = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").ClearFormats()
$ws.Range("D49").Formula = "'8.962"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.06%  "

$ws.Range("D50").Formula = "'33.81"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Formula = "'0.05667"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.85%  "
